$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------------------
# The sheet currently has 80 rows of inline text in column A (A1 header + 79 data rows,
# A2 being a stray "Executive Summary / INTRODUCTION" row, and the long air-pollution
# paragraph split across A15/A16). The target layout is 78 rows (A1 header + 77 data rows)
# split across two columns: A = a 1-based sequence number (as text), B = the paragraph text.
# ----------------------------------------------------------------------------------

# Capture the two halves of the split air-pollution paragraph (rows 15 & 16) up front, while
# the original row numbering is still in effect, and glue them back into one paragraph.
$airPart1 = $ws.Cells.Item(15, 1).Text
$airPart2 = $ws.Cells.Item(16, 1).Text
$airMerged = $airPart1 + "012. " + $airPart2

# Drop the now-redundant second half (row 16), then drop the stray "Executive Summary" row (row 2).
# Row 16 is removed first so that removing row 2 (which shifts everything below it up by one) does
# not change row 16's index out from under us.
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(2).Delete()

# After both deletions: old row 15 (first half of the air-pollution paragraph) now sits at row 14,
# and the 77 remaining paragraphs occupy rows 2-78 of column A.

# Give B1 the same header style as A1 (bold/centered/bordered), then set the new header captions.
$ws.Range("A1").Copy($ws.Range("B1"))
$ws.Range("A1").Value = 'Number'
$ws.Range("B1").Value = 'Text'

# Walk the 77 content rows: stash each paragraph (substituting the glued-together air-pollution
# paragraph at row 14), then overwrite column A with a 1-based sequence number stored as TEXT
# (matching the source data - "1", "2", "3", ... not numeric 1, 2, 3). The temporary "@" (Text)
# number format forces the numeral strings to stick as text instead of auto-converting to
# numbers; ClearFormats drops that temporary formatting again once the values are in place.
# This full pass over column A has to finish *before* column B is touched below, otherwise the
# row auto-height recalculation triggered by the long column-B text interacts badly with
# ClearFormats and leaves stray formatting behind.
$paragraphs = @()
$ws.Range("A2:A78").NumberFormat = "@"
for ($i = 2; $i -le 78; $i++) {
    if ($i -eq 14) {
        $paragraphs += $airMerged
    } else {
        $paragraphs += $ws.Cells.Item($i, 1).Text
    }
    $ws.Cells.Item($i, 1).Value = [string]($i - 1)
}
$ws.Range("A2:A78").ClearFormats()

# Finally, move the stashed paragraph text into column B.
for ($i = 2; $i -le 78; $i++) {
    $ws.Cells.Item($i, 2).Value = $paragraphs[$i - 2]
}
